$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 1
